# Update '想去人数' (F column) figures per the latest site crawl (commit 456a3b4)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2664
$ws.Range("F5").Value = 1488
$ws.Range("F8").Value = 535
$ws.Range("F13").Value = 9085
$ws.Range("F14").Value = 392
$ws.Range("F15").Value = 2499
$ws.Range("F18").Value = 178
$ws.Range("F23").Value = 1000
$ws.Range("F24").Value = 2083
$ws.Range("F25").Value = 2167
$ws.Range("F27").Value = 1874
$ws.Range("F29").Value = 1926
$ws.Range("F31").Value = 1178
$ws.Range("F32").Value = 269
$ws.Range("F33").Value = 146
$ws.Range("F34").Value = 204
$ws.Range("F38").Value = 288
$ws.Range("F39").Value = 479
$ws.Range("F41").Value = 19
$ws.Range("F42").Value = 227
$ws.Range("F43").Value = 1369
$ws.Range("F44").Value = 292
$ws.Range("F45").Value = 4
$ws.Range("F46").Value = 4
$ws.Range("F47").Value = 607

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 20
$ws.Range("F5").Value = 19

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2664
$ws.Range("F5").Value = 1488
$ws.Range("F9").Value = 535
$ws.Range("F12").Value = 9085
$ws.Range("F13").Value = 392
$ws.Range("F14").Value = 2499
$ws.Range("F15").Value = 20
$ws.Range("F19").Value = 178
$ws.Range("F23").Value = 1000
$ws.Range("F24").Value = 2167
$ws.Range("F25").Value = 1874
$ws.Range("F26").Value = 1926
$ws.Range("F28").Value = 1179
$ws.Range("F29").Value = 269
$ws.Range("F30").Value = 146
$ws.Range("F31").Value = 204
$ws.Range("F35").Value = 288
$ws.Range("F36").Value = 479
$ws.Range("F37").Value = 19
$ws.Range("F41").Value = 19
$ws.Range("F42").Value = 227
$ws.Range("F44").Value = 1369
$ws.Range("F46").Value = 292
$ws.Range("F47").Value = 4
$ws.Range("F48").Value = 607
